# Weekly refresh of the Fruta/Hortaliza "Coco" subset: the per-row
# Fecha/Calidad/Volumen/Precio.../Origen values get reshuffled among the
# existing data rows (2..41) while the descriptive columns (Mercado,
# Región, Codreg, Tipo, Producto, Categoría, Variedad, Unidad, Kg/unidad)
# stay put on each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: new row -> source row (old data) that feeds its shuffled columns.
$rowMap = @{
    2 = 6;  3 = 14; 4 = 10; 5 = 18; 6 = 41; 7 = 8;  8 = 32; 9 = 27; 10 = 19;
    11 = 5; 12 = 2; 13 = 37; 14 = 40; 15 = 25; 16 = 20; 17 = 34; 18 = 28;
    19 = 9; 20 = 33; 21 = 4; 22 = 29; 23 = 21; 24 = 3; 25 = 11; 26 = 16;
    27 = 24; 28 = 31; 29 = 26; 30 = 30; 31 = 7; 32 = 39; 33 = 36; 34 = 22;
    35 = 17; 36 = 35; 37 = 12; 38 = 15; 39 = 23; 40 = 38; 41 = 13
}

# Columns that travel together as a group for each logical record.
# D=Fecha, L=Calidad, M=Volumen, N=Precio minimo, O=Precio maximo,
# P=Precio promedio ponderado, R=Origen, S=Precio $/Kg
$cols = @(4, 12, 13, 14, 15, 16, 18, 19)

# Snapshot the original values for every row/column we might read from,
# since the permutation has multi-row cycles and must not read values
# that were already overwritten.
$snapshot = @{}
for ($r = 2; $r -le 41; $r++) {
    foreach ($c in $cols) {
        $snapshot["$r-$c"] = $ws.Cells.Item($r, $c).Value2
    }
}

foreach ($r in 2..41) {
    $src = $rowMap[$r]
    foreach ($c in $cols) {
        $ws.Cells.Item($r, $c).Value2 = $snapshot["$src-$c"]
    }
}
